# Mapa.xlsx — "mapa de monstruos" edit
#
# Updates the "idsm" (monster-seed id) column for the three location rows
# that previously had no seed assigned (0), and moves the active selection
# on Hoja1 from D5 to D10.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Column C ("idsm") holds the id of the monster seed that spawns at each
# location. Rows 3-5 were placeholders (0); assign the new seed ids.
$ws.Range("C3").Value = 1
$ws.Range("C4").Value = 3
$ws.Range("C5").Value = 5

# Move the current selection to D10, as left by the author after the edit.
$ws.Activate()
$ws.Range("D10").Select()
